$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Rows 153 and 154 had their match data (columns F:V) swapped.
#    Columns A:E (Indice, pais, torneio, temporada, data_partida) stay put.
#    Capture the current F:V values/ (both match-detail columns) before
#    overwriting anything, then write them back swapped.
# ------------------------------------------------------------------
$row153 = @{}
$row154 = @{}
for ($col = 6; $col -le 22; $col++) {
    $row153[$col] = $ws.Cells.Item(153, $col).Value2
    $row154[$col] = $ws.Cells.Item(154, $col).Value2
}

for ($col = 6; $col -le 22; $col++) {
    $ws.Cells.Item(153, $col).Value = $row154[$col]
    $ws.Cells.Item(154, $col).Value = $row153[$col]
}

# ------------------------------------------------------------------
# 2) Append a new match as row 159 (Indice 158), copying the
#    formatting of the last existing data row (158) so the new row's
#    styles (bold/bordered index cell, date-formatted data_partida
#    cell) match the rest of the table.
# ------------------------------------------------------------------
$ws.Range("A158:V158").Copy()
$ws.Range("A159:V159").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(159, 1).Value = 158
$ws.Cells.Item(159, 2).Value = "poland"
$ws.Cells.Item(159, 3).Value = "division-2"
$ws.Cells.Item(159, 4).Value = "2023-2024"
$ws.Cells.Item(159, 5).Value = 45250.66666666666
$ws.Cells.Item(159, 6).Value = "LKS Lodz II"
$ws.Cells.Item(159, 7).Value = 1
$ws.Cells.Item(159, 8).Value = "Pogon Siedlce"
$ws.Cells.Item(159, 9).Value = 2
$ws.Cells.Item(159, 10).Value = 2.39
$ws.Cells.Item(159, 11).Value = "19/11/2023 04:12"
$ws.Cells.Item(159, 12).Value = 2.63
$ws.Cells.Item(159, 13).Value = "20/11/2023 15:58"
$ws.Cells.Item(159, 14).Value = 3.31
$ws.Cells.Item(159, 15).Value = "19/11/2023 04:12"
$ws.Cells.Item(159, 16).Value = 3.56
$ws.Cells.Item(159, 17).Value = "20/11/2023 15:50"
$ws.Cells.Item(159, 18).Value = 2.54
$ws.Cells.Item(159, 19).Value = "19/11/2023 04:12"
$ws.Cells.Item(159, 20).Value = 2.44
$ws.Cells.Item(159, 21).Value = "20/11/2023 15:58"
$ws.Cells.Item(159, 22).Value = "https://www.betexplorer.com/football/poland/division-2/lks-lodz-pogon-siedlce/xChRnjB8/"
